# The author renamed the internal "name" attribute used for the three
# logo pictures that live in the document's header/footer parts:
#   - footer (docPr id="2")  PearsonLogo: image2.png -> image1.png
#   - footer (docPr id="3")  PearsonLogo: image2.png -> image1.png
#   - header (docPr id="1")  BTec_Logo-Orange: image1.jpg -> image2.jpg
#
# InlineShape has no writable .Name in the Word object model, so we
# round-trip the owning header/footer Range's WordOpenXML and patch the
# "name" attribute on the wp:docPr / pic:cNvPr pair for that picture.

$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlineImage($range, [string]$oldName, [string]$newName) {
    $xml = $range.WordOpenXML
    $patched = $xml.Replace('name="' + $oldName + '"', 'name="' + $newName + '"')
    if ($patched -ne $xml) {
        $range.WordOpenXML = $patched
    }
}

# Footer containing docPr id="2" (PearsonLogo) -> image1.png
Rename-InlineImage $sec.Footers(1).Range "image2.png" "image1.png"

# Footer containing docPr id="3" (PearsonLogo) -> image1.png
Rename-InlineImage $sec.Footers(2).Range "image2.png" "image1.png"

# Header containing docPr id="1" (BTec_Logo-Orange) -> image2.jpg
Rename-InlineImage $sec.Headers(2).Range "image1.jpg" "image2.jpg"
